# Edit slide 3 ("DATASET AND OERVIEW") of the presentation:
#  - reposition/resize the Title and Content placeholders
#  - rewrite the content placeholder's text with new formatting (bold /
#    red highlights) and two new paragraphs with hyperlinks to the
#    dataset source files.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---------------------------------------------------------------------
# 1. Title placeholder (shape id=2): move up & slightly resize
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Top    = 18.509213598425198   # 235067 EMU
$title.Width  = 624.3497627795275    # 7929242 EMU
$title.Height = 81.38614273228346    # 1033604 EMU

# ---------------------------------------------------------------------
# 2. Content placeholder (shape id=3): move up & enlarge, then
#    rebuild the text with the new runs/formatting.
# ---------------------------------------------------------------------
$content = $s.Shapes.Item(2)
$content.Top    = 99.92173428346456  # 1269006 EMU
$content.Width  = 698.4787301574803  # 8870680 EMU
$content.Height = 410.13929833858265 # 5208769 EMU

$tr = $content.TextFrame.TextRange

$para1a = "The dataset has taken from "
$para1b = "The Canadian Institute for Cybersecurity"
$para1c = " (CIC) generate the dataset in collaboration with "
$para1d = "The "
$para1e = "Communications Security Establishment"
$para1f = " (CSE) of network events, some benign and other malicious. "

$para2 = "By using CICFlowMeter-V3, around 80 network traffic features were extracted from each attack."

$para3 = "The goal of this project is to use this dataset to identify cybersecurity threats and furthermore to classify threat with as much accuracy as possible using machine learning"

$githubUrl = "https://github.com/AsimGull/Data-Science-Projects/blob/main/Machine%20learning/Natural%20language%20processing/DDoS%20Attack%20Classification/data/DDos.csv"

$origDatasetLabel = "Original dataset"

$unbUrl = "https://www.unb.ca/cic/datasets/index.html"

$para1 = $para1a + $para1b + $para1c + $para1d + $para1e + $para1f

$fullText = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $githubUrl + "`r" + "`r" + $origDatasetLabel + "`r" + $unbUrl + "`r" + "`r"

$tr.Text = $fullText

# --- paragraph 1: apply bold / red highlighting to two phrases -------
$pos = 1
$r1a = $tr.Characters($pos, $para1a.Length); $pos += $para1a.Length
$r1b = $tr.Characters($pos, $para1b.Length); $pos += $para1b.Length
$r1c = $tr.Characters($pos, $para1c.Length); $pos += $para1c.Length
$r1d = $tr.Characters($pos, $para1d.Length); $pos += $para1d.Length
$r1e = $tr.Characters($pos, $para1e.Length); $pos += $para1e.Length
$r1f = $tr.Characters($pos, $para1f.Length); $pos += $para1f.Length

$r1b.Font.Bold = $true
$r1b.Font.Color.RGB = 255          # FF0000 (BGR-packed red)

$r1d.Font.Bold = $true

$r1e.Font.Bold = $true
$r1e.Font.Color.RGB = 255          # FF0000

# move past the paragraph mark
$pos += 1

# --- paragraph 2 & 3: plain text, nothing extra to format ------------
$pos += $para2.Length
$pos += 1
$pos += $para3.Length
$pos += 1

# --- paragraph 4: github URL, hyperlinked ------------------------------
$rGithub = $tr.Characters($pos, $githubUrl.Length)
$rGithub.ActionSettings.Item(1).Hyperlink.Address = $githubUrl
$pos += $githubUrl.Length
$pos += 1

# --- paragraph 5: blank ------------------------------------------------
$pos += 1

# --- paragraph 6: "Original dataset" bold + underline -----------------
$rOrig = $tr.Characters($pos, $origDatasetLabel.Length)
$rOrig.Font.Bold = $true
$rOrig.Font.Underline = $true
$pos += $origDatasetLabel.Length
$pos += 1

# --- paragraph 7: unb.ca URL, hyperlinked ------------------------------
$rUnb = $tr.Characters($pos, $unbUrl.Length)
$rUnb.ActionSettings.Item(1).Hyperlink.Address = $unbUrl
$pos += $unbUrl.Length
$pos += 1

# --- paragraphs 8 & 9: blank -------------------------------------------

# --- remove bullets from paragraphs 4-9 (the link/label block) --------
for ($i = 4; $i -le 9; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Visible = $false
}

Write-Output "slide 3 updated"
